# Weekly crypto price/volume refresh (GitHub Actions data pull).
# D-column values are plain numeric-looking strings (e.g. "184.90", "3.493.82")
# that must stay TEXT, matching the source sheet's inlineStr cells; a leading
# apostrophe is Excel's standard quote-prefix trick to force text and preserve
# exact formatting (trailing zeros, dotted thousand separators, etc.).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '''67.415.75'
$ws.Range('E2').Value = '  +4.87%  '

# Row 3
$ws.Range('D3').Value = '''3.486.07'
$ws.Range('E3').Value = '  +4.79%  '

# Row 4
$ws.Range('E4').Value = '  +0.19%  '

# Row 5
$ws.Range('D5').Value = '''588.11'

# Row 6
$ws.Range('D6').Value = '''184.90'
$ws.Range('E6').Value = '  +6.66%  '

# Row 7
$ws.Range('D7').Value = '''0.634'
$ws.Range('E7').Value = '  +2.24%  '

# Row 8
$ws.Range('D8').Value = '''3.485.83'
$ws.Range('E8').Value = '  +5.01%  '

# Row 9
$ws.Range('E9').Value = '  +0.13%  '

# Row 10
$ws.Range('D10').Value = '''0.174'
$ws.Range('E10').Value = '  +2.44%  '

# Row 11
$ws.Range('D11').Value = '''0.651'
$ws.Range('E11').Value = '  +3.32%  '

# Row 12
$ws.Range('D12').Value = '''56.18'
$ws.Range('E12').Value = '  +5.35%  '

# Row 13
$ws.Range('D13').Value = '''0.0000280'
$ws.Range('E13').Value = '  +0.37%  '

# Row 14
$ws.Range('D14').Value = '''9.46'
$ws.Range('E14').Value = '  +4.20%  '

# Row 15
$ws.Range('D15').Value = '''4.040.90'
$ws.Range('E15').Value = '  +4.98%  '

# Row 16
$ws.Range('D16').Value = '''18.84'
$ws.Range('E16').Value = '  +4.10%  '

# Row 17
$ws.Range('D17').Value = '''3.492.69'
$ws.Range('E17').Value = '  +5.28%  '

# Row 18
$ws.Range('D18').Value = '''67.477.91'
$ws.Range('E18').Value = '  +4.88%  '

# Row 19
$ws.Range('D19').Value = '''12.18'
$ws.Range('E19').Value = '  +4.10%  '

# Row 20
$ws.Range('D20').Value = '''0.119'
$ws.Range('E20').Value = '  -1.22%  '

# Row 21
$ws.Range('E21').Value = '  +3.67%  '

# Row 22
$ws.Range('D22').Value = '''490.52'
$ws.Range('E22').Value = '  +8.06%  '

# Row 23
$ws.Range('D23').Value = '''5.45'
$ws.Range('E23').Value = '  +6.17%  '

# Row 24
$ws.Range('D24').Value = '''16.67'
$ws.Range('E24').Value = '  +18.96%  '

# Row 25
$ws.Range('D25').Value = '''4.48'
$ws.Range('E25').Value = '  +10.41%  '

# Row 26
$ws.Range('D26').Value = '''90.16'
$ws.Range('E26').Value = '  +3.35%  '

# Row 27
$ws.Range('D27').Value = '''2.96'
$ws.Range('E27').Value = '  +3.37%  '

# Row 28
$ws.Range('D28').Value = '''10.99'
$ws.Range('E28').Value = '  +3.50%  '

# Row 29
$ws.Range('D29').Value = '''9.15'
$ws.Range('E29').Value = '  +6.79%  '

# Row 30
$ws.Range('D30').Value = '''31.57'
$ws.Range('E30').Value = '  +2.12%  '

# Row 31
$ws.Range('D31').Value = '''7.15'
$ws.Range('E31').Value = '  +9.81%  '

# Row 32
$ws.Range('B32').Value = 'OKB'
$ws.Range('C32').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D32').Value = '''64.48'
$ws.Range('E32').Value = '  +4.77%  '

# Row 33
$ws.Range('B33').Value = 'Cosmos'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D33').Value = '''11.75'
$ws.Range('E33').Value = '  +3.07%  '

# Row 34
$ws.Range('D34').Value = '''595.08'
$ws.Range('E34').Value = '  +5.34%  '

# Row 35
$ws.Range('D35').Value = '''0.111'
$ws.Range('E35').Value = '  +4.48%  '

# Row 36
$ws.Range('E36').Value = '  +6.62%  '

# Row 37
$ws.Range('E37').Value = '  -0.12%  '

# Row 38
$ws.Range('D38').Value = '''36.63'
$ws.Range('E38').Value = '  +3.86%  '

# Row 39
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').Value = '''3.56'
$ws.Range('E39').Value = '  +1.64%  '

# Row 40
$ws.Range('B40').Value = 'TheGraph'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D40').Value = '''0.386'
$ws.Range('E40').Value = '  +5.84%  '

# Row 41
$ws.Range('D41').Value = '''0.0₃0768'
$ws.Range('E41').Value = '  +5.63%  '

# Row 42
$ws.Range('D42').Value = '''3.252.51'
$ws.Range('E42').Value = '  +6.57%  '

# Row 43
$ws.Range('D43').Value = '''2.91'
$ws.Range('E43').Value = '  +6.31%  '

# Row 44
$ws.Range('E44').Value = '  +3.45%  '

# Row 45
$ws.Range('D45').Value = '''2.53'
$ws.Range('E45').Value = '  +3.71%  '

# Row 46
$ws.Range('D46').Value = '''3.26'
$ws.Range('E46').Value = '  +1.88%  '

# Row 47
$ws.Range('D47').Value = '''2.75'
$ws.Range('E47').Value = '  +21.91%  '

# Row 48
$ws.Range('D48').Value = '''0.136'
$ws.Range('E48').Value = '  +1.72%  '

# Row 49
$ws.Range('D49').Value = '''3.26'
$ws.Range('E49').Value = '  +12.36%  '

# Row 50
$ws.Range('E50').Value = '  +7.78%  '

# Row 51
$ws.Range('D51').Value = '''0.999'
$ws.Range('E51').Value = '  -0.06%  '
